$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "301.55"
Set-TextValue "E2" "-0.71%"
Set-TextValue "D3" "31.43"
Set-TextValue "E3" "-2.11%"
Set-TextValue "D4" "5.092"
Set-TextValue "E4" "-3.23%"
Set-TextValue "D5" "0.07367"
Set-TextValue "E5" "-2.16%"
Set-TextValue "D6" "2.269"
Set-TextValue "E6" "51.48%"
Set-TextValue "D7" "7.960"
Set-TextValue "E7" "0.71%"
Set-TextValue "D8" "3.782"
Set-TextValue "E8" "-0.85%"
Set-TextValue "D9" "0.9183"
Set-TextValue "E9" "-0.36%"
Set-TextValue "D10" "0.1708"
Set-TextValue "E10" "0.35%"
Set-TextValue "D11" "0.07547"
Set-TextValue "E11" "-5.46%"
Set-TextValue "D12" "0.08063"
Set-TextValue "E12" "0.46%"
Set-TextValue "D13" "0.03022"
Set-TextValue "E13" "-0.47%"
Set-TextValue "D14" "0.09930"
Set-TextValue "E14" "0.26%"
Set-TextValue "D15" "0.001498"
Set-TextValue "E15" "-1.03%"
Set-TextValue "D16" "0.006070"
Set-TextValue "E16" "-7.10%"
Set-TextValue "E17" "-0.03%"
Set-TextValue "D18" "2.224"
Set-TextValue "E18" "-0.42%"
Set-TextValue "D19" "0.3306"
Set-TextValue "E19" "0.21%"
Set-TextValue "E20" "0.05%"
Set-TextValue "D21" "4.641"
Set-TextValue "E21" "3.70%"
Set-TextValue "D22" "0.04653"
Set-TextValue "E22" "1.24%"
Set-TextValue "D23" "0.1565"
Set-TextValue "E23" "-3.28%"
Set-TextValue "D24" "0.001224"
Set-TextValue "E24" "0.66%"
Set-TextValue "D25" "0.004493"
Set-TextValue "E25" "0.75%"
Set-TextValue "D26" "0.0001298"
Set-TextValue "E26" "-7.13%"
Set-TextValue "E27" "50.43%"
Set-TextValue "D39" "0.01747"
Set-TextValue "E39" "2.60%"
Set-TextValue "D40" "0.04509"
Set-TextValue "E40" "0.44%"
Set-TextValue "D41" "0.007222"
Set-TextValue "E41" "3.75%"
Set-TextValue "D42" "0.1349"
Set-TextValue "E42" "-0.26%"
Set-TextValue "D43" "0.002226"
Set-TextValue "E43" "4.22%"
Set-TextValue "D45" "0.00006287"
Set-TextValue "E45" "1.93%"
Set-TextValue "E46" "15.83%"
